# CIERRE 13 Nov 2021
# Updates the weekly payroll figures on "Hoja1" and moves the saved
# cursor/scroll position to reflect where the author finished editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Block 1 (rows 2-15): "Extra " amount for the first payslip ---
$ws.Range("K4").Value = 1300

# --- Block 2 (rows 20-26): "DESCUENTO" amount for the second payslip ---
$ws.Range("K21").Value = 3080

# --- Block 3 (rows 36-41): figures for the third payslip ---
$ws.Range("E38").Value = 2200
# The PRESTAMO label/amount on row 39 is cleared out (no loan this week).
$ws.Range("J39").ClearContents()
$ws.Range("K39").Value = 0

# --- Move the saved selection/cursor to where editing finished ---
$ws.Range("G44").Select()

Write-Output "done"
